$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$s71 = @'
<p>In this informational e-book, Megan M. Gunderson describes the pivotal plans and actions of President John F. Kennedy during his administration, focusing on domestic initiatives, international relations, and the Cold War context. Through detailed accounts of events such as the Bay of Pigs invasion and the Cuban Missile Crisis, the author conveys Kennedy's strategic decisions and their implications.</p>
'@

$s72 = @'
<ul>
<li>poverty (noun): the state of being extremely poor.</li>
<li>civil rights (noun): the rights of citizens to political and social freedom and equality.</li>
<li>Cold War (proper noun): the period of geopolitical tension between the Soviet Union and the United States after World War II.</li>
<li>Soviet Union (proper noun): a former federation of communist republics in Eastern Europe and northern Asia.</li>
<li>Communist Soviet Union (proper noun): a reference to the Soviet Union as a communist state.</li>
<li>relationships (noun): the way in which two or more people or groups regard and behave toward each other.</li>
<li>Communism (proper noun): a political and economic ideology advocating for a classless system in which the means of production are communally owned.</li>
<li>overthrow (verb): to remove forcibly from power.</li>
<li>Communist (proper adjective): relating to or advocating for communism.</li>
<li>Bay of Pigs (proper noun): a coastal area in Cuba known for the failed invasion of Cubans in 1961.</li>
<li>Communist China (proper noun): a reference to the People's Republic of China as a communist state.</li>
<li>weapons (noun): objects used to inflict physical harm or damage.</li>
<li>nuclear (adjective): relating to the nucleus of an atom, often associated with weapons of mass destruction.</li>
<li>supplies (noun): the necessary items required for a particular purpose.</li>
<li>advisers (noun): people who provide guidance or recommendations.</li>
<li>showdown (noun): a decisive confrontation or contest.</li>
<li>Cuban Missile Crisis (proper noun): a 1962 confrontation between the United States and the Soviet Union over missile deployment in Cuba.</li>
</ul>
'@

$s73 = @'
<p>Making Plans (excerpt from John F. Kennedy) <br>By Megan M. Gunderson, Abdo Publishing <br><br>John F. Kennedy became President of the United States on January 20, 1961. He had many plans when he took office. These included programs for higher wages and aid for people experiencing *poverty*. He also hoped to improve *civil rights*. <br><br>During this time, America was part of the *Cold War* with the *Soviet Union*. Another plan was to put a person on the moon before the *Communist Soviet Union* did. As a result, Kennedy greatly improved the US space program. <br><br>*Relationships* with other countries were another concern for Kennedy. American officials wanted to stop the spread of *Communism*. So, they trained a group of Cubans to *overthrow* Cuba&rsquo;s *Communist* leader. In 1961, the group invaded Cuba at the *Bay of Pigs*. However, the effort failed. <br><br>Then, in 1962, *Communist* China invaded India. Kennedy sent *weapons* to India&rsquo;s army. Meanwhile, South Vietnam was fighting *Communist* North Vietnam. Kennedy sent thousands of US military *advisers* to help the South Vietnamese. <br><br>In October 1962, Kennedy learned *nuclear* missile bases were being built in Cuba. He wanted to stop Soviet ships from bringing *supplies* there. So, Kennedy stated that the US Navy would block Cuba&rsquo;s coast. <br><br>The Soviet leader said his country would guard its shipping rights. A *nuclear* war could begin. But Kennedy held his ground. The Soviet ships stayed away from the blocked area. Kennedy won the *showdown*! This event became known as the *Cuban Missile Crisis*.</p>
'@

$s74 = @'
<p>1. What was one of John F. Kennedy&rsquo;s key plans regarding the space race during the Cold War?<br>A. To increase military presence in Europe<br>B. To send a person to the moon before the Soviet Union<br>C. To strengthen trade relations with Communist countries<br>D. To establish a United Nations peacekeeping mission<br>[CCSS.RI.5.1; TEKS.ELAR.5.7(C)]<br><br>2. How did the Cold War influence U.S. foreign policy in Latin America?<br>A. The U.S. sought to promote trade with Communist nations<br>B. The U.S. supported a Cuban invasion to help Communist leaders<br>C. The U.S. trained Cubans to overthrow their Communist leader<br>D. The U.S. engaged in direct military intervention in all Latin American countries<br>[CCSS.RI.4.1; TEKS.ELAR.4.7(C)]<br><br>3. What was a major international event during the Cold War that escalated tensions between the United States and the Soviet Union?<br>A. The invasion of Europe in World War II<br>B. The Cuban Missile Crisis<br>C. The founding of the United Nations<br>D. The signing of the Treaty of Versailles<br>[CCSS.RI.3.1; TEKS.ELAR.3.7(C)]<br><br>4. In response to Communist aggression in other countries, what action did Kennedy take concerning South Vietnam?<br>A. He withdrew all military advisers<br>B. He sent vaccines and aid to the population<br>C. He sent thousands of U.S. military advisers to support the South<br>D. He called for a cease-fire between North and South Vietnam<br>[BEST.ELA.K12.EE.3.1; CCSS.RI.3.1]<br><br>5. What was the result of Kennedy's decision to block Cuba's coast during the Cuban Missile Crisis?<br>A. A nuclear war broke out between the U.S. and the Soviet Union<br>B. The Soviet ships stayed away from the blocked area<br>C. The United States launched missiles at Cuba<br>D. Kennedy changed his policy on Communism<br>[CCSS.RI.4.1; TEKS.ELAR.4.7(C)]<br><br></p>
'@

$s75 = @'
<p>In this informational e-book, Megan M. Gunderson describes the plans and challenges faced by John F. Kennedy during his presidency through factual recounting and analysis of key historical events, including his domestic policies and foreign relations strategies.</p>
'@

$s76 = @'
<ul>
<li>poverty (noun): the state of being extremely poor.</li>
<li>civil rights (noun): the rights of citizens to political and social freedom and equality.</li>
<li>Cold War (proper noun): the period of geopolitical tension between the Soviet Union and the United States and their respective allies after World War II.</li>
<li>Soviet Union (proper noun): a former federal socialist state in eastern Europe and northern Asia that existed from 1922 to 1991.</li>
<li>Communist (proper adjective): relating to or characterized by communism, a political and economic ideology.</li>
<li>space (noun): the vast, seemingly infinite expanse that exists beyond Earth and between celestial bodies.</li>
<li>relationships (noun): the way in which two or more people or groups regard and behave toward each other.</li>
<li>Communism (proper noun): a political and economic ideology advocating for a classless system in which the means of production are owned communally.</li>
<li>Bay of Pigs (proper noun): a failed military invasion of Cuba undertaken by the United States in 1961.</li>
<li>Communist (proper adjective): relating to or characterized by communism, a political and economic ideology.</li>
<li>weapons (noun): tools or devices used for fighting or attacking.</li>
<li>Communist (proper adjective): relating to or characterized by communism, a political and economic ideology.</li>
<li>nuclear missile (noun): a missile capable of delivering a nuclear warhead to a target.</li>
<li>supplies (noun): the provisions or materials required for a particular purpose.</li>
<li>shipping (noun): the act of transporting goods or cargo.</li>
<li>nuclear (adjective): relating to the atomic nucleus or nuclear energy.</li>
<li>showdown (noun): a decisive confrontation or contest.</li>
<li>Cuban Missile Crisis (proper noun): a 13-day confrontation in October 1962 between the United States and the Soviet Union concerning Soviet ballistic missiles deployed in Cuba.</li>
</ul>
'@

$s77 = @'
<p>Making Plans (excerpt from John F. Kennedy) <br>By Megan M. Gunderson, Abdo Publishing <br><br>John F. Kennedy became President of the United States on January 20, 1961. He had many plans when he took office. These included programs for higher wages and aid for people experiencing *poverty*. He also hoped to improve *civil rights*. <br><br>During this time, America was part of the *Cold War* with the *Soviet Union*. Another plan was to put a person on the moon before the *Communist* Soviet Union did. As a result, Kennedy greatly improved the US *space* program. <br><br>*Relationships* with other countries were another concern for Kennedy. American officials wanted to stop the spread of *Communism*. So, they trained a group of Cubans to overthrow Cuba&rsquo;s *Communist* leader. In 1961, the group invaded Cuba at the *Bay of Pigs*. However, the effort failed. <br><br>Then, in 1962, *Communist* China invaded India. Kennedy sent *weapons* to India&rsquo;s army. Meanwhile, South Vietnam was fighting *Communist* North Vietnam. Kennedy sent thousands of US military *advisers* to help the South Vietnamese. <br><br>In October 1962, Kennedy learned *nuclear missile* bases were being built in Cuba. He wanted to stop Soviet ships from bringing *supplies* there. So, Kennedy stated that the US Navy would block Cuba&rsquo;s coast. <br><br>The Soviet leader said his country would guard its *shipping* rights. A *nuclear* war could begin. But Kennedy held his ground. The Soviet ships stayed away from the blocked area. Kennedy won the *showdown*! This event became known as the *Cuban Missile Crisis*.</p>
'@

$s78 = @'
<p>1. How did John F. Kennedy's actions during the Cold War reflect America's desire to contain Communism?<br>A. By improving civil rights programs <br>B. By sending weapons to allies fighting Communism <br>C. By cooperating with the Soviet Union <br>D. By focusing solely on domestic issues <br>[CCSS.RI.4.1; TEKS.ELAR.4.7(C)] <br><br>2. What was one consequence of the Bay of Pigs invasion led by Kennedy's administration?<br>A. Strengthened U.S. relations with Cuba <br>B. Increased support for Communism in Latin America <br>C. Successful overthrow of the Cuban government <br>D. Heightened tensions between the U.S. and the Soviet Union <br>[CCSS.RI.3.1; TEKS.ELAR.3.7(C)] <br><br>3. How did Kennedy's response to the Cuban Missile Crisis exemplify a pivotal moment in U.S.-Soviet relations?<br>A. It demonstrated American willingness to retreat during conflict <br>B. It showed the effectiveness of military intervention in foreign crises <br>C. It marked a successful negotiation that avoided nuclear war <br>D. It led to an immediate partnership with the Soviet Union <br>[CCSS.RI.5.1; TEKS.ELAR.5.7(C)] <br><br>4. What role did the space race play in shaping international relations during the Cold War?<br>A. It caused nations to collaborate on technology <br>B. It increased competition and national prestige between the U.S. and the Soviet Union <br>C. It led to combined military ventures <br>D. It reduced tensions between ideologies <br>[CCSS.RI.3.1; BEST.ELA.K12.EE.3.1] <br><br>5. What inference can be made about Kennedy's military involvement in Vietnam and its impact on international relations?<br>A. It isolated the U.S. from global conflicts <br>B. It encouraged stronger alliances in Southeast Asia <br>C. It showed a commitment to fighting against communism globally <br>D. It indicated a lack of interest in foreign military engagements <br>[CCSS.RI.4.1; TEKS.ELAR.4.7(C)] <br><br></p>
'@

$s79 = @'
<p>In this informational e-book, Megan M. Gunderson describes the key plans and actions of President John F. Kennedy during his administration, highlighting his initiatives in civil rights, the space program, and international relations through concise historical accounts and factual details.</p>
'@

$s80 = @'
<ul>
<li>poverty (noun): the state of being extremely poor.</li>
<li>civil rights (noun): the rights of citizens to political and social freedom and equality.</li>
<li>Cold War (proper noun): a period of geopolitical tension between the Soviet Union and the United States after World War II.</li>
<li>Soviet Union (proper noun): a former federal socialist state in northern Eurasia that existed from 1922 to 1991.</li>
<li>Communist (proper adjective): relating to the political and economic ideology of communism.</li>
<li>space program (noun): a program for the exploration of outer space.</li>
<li>Communism (proper noun): a political and economic ideology advocating for class war and leading to a society in which all property is publicly owned.</li>
<li>overthrow (verb): to remove from power, usually by force.</li>
<li>Bay of Pigs (proper noun): a failed invasion of Cuba in 1961 by Cuban exiles, supported by the U.S. government.</li>
<li>Communist (proper adjective): relating to or characteristic of communism.</li>
<li>weapons (noun): tools or instruments used for combat or warfare.</li>
<li>advisers (noun): individuals who provide advice, typically in a professional context.</li>
<li>nuclear (adjective): relating to the nucleus of an atom, often used in the context of nuclear energy or weapons.</li>
<li>showdown (noun): a decisive confrontation or contest.</li>
<li>Cuban Missile Crisis (proper noun): a 13-day confrontation in October 1962 between the United States and the Soviet Union over nuclear missiles deployed in Cuba.</li>
</ul>
'@

$s81 = @'
<p>Making Plans (excerpt from John F. Kennedy) <br>By Megan M. Gunderson, Abdo Publishing <br><br>John F. Kennedy became President of the United States on January 20, 1961. He had many plans when he took office. These included programs for higher wages and aid for people experiencing *poverty*. He also hoped to improve *civil rights*. <br><br>During this time, America was part of the *Cold War* with the *Soviet Union*. Another plan was to put a person on the moon before the *Communist* Soviet Union did. As a result, Kennedy greatly improved the US *space program*. <br><br>Relationships with other countries were another concern for Kennedy. American officials wanted to stop the spread of *Communism*. So, they trained a group of Cubans to *overthrow* Cuba&rsquo;s Communist leader. In 1961, the group invaded Cuba at the *Bay of Pigs*. However, the effort failed. <br><br>Then, in 1962, *Communist* China invaded India. Kennedy sent *weapons* to India&rsquo;s army. Meanwhile, South Vietnam was fighting Communist North Vietnam. Kennedy sent thousands of US military *advisers* to help the South Vietnamese. <br><br>In October 1962, Kennedy learned *nuclear* missile bases were being built in Cuba. He wanted to stop Soviet ships from bringing supplies there. So, Kennedy stated that the US Navy would block Cuba&rsquo;s coast. <br><br>The Soviet leader said his country would guard its shipping rights. A nuclear war could begin. But Kennedy held his ground. The Soviet ships stayed away from the blocked area. Kennedy won the *showdown*! This event became known as the *Cuban Missile Crisis*.</p>
'@

$s82 = @'
<p>1. What was one of John F. Kennedy's priorities when he became President in 1961? <br>A. To reduce military presence in Vietnam <br>B. To improve civil rights <br>C. To decrease space exploration <br>D. To strengthen ties with the Soviet Union <br>[CCSS.RI.3.1; TEKS.ELAR.3.7(C); BEST.ELA.K12.EE.3.1; CCSS.RI.4.1; TEKS.ELAR.4.7(C); CCSS.RI.5.1; TEKS.ELAR.5.7(C)] <br><br>2. How did Kennedy respond to the spread of Communism in Cuba? <br>A. By establishing a peaceful negotiation <br>B. By training a group of Cubans to invade Cuba <br>C. By withdrawing all American forces <br>D. By forming an alliance with the Soviets <br>[CCSS.RI.3.1; TEKS.ELAR.3.7(C); BEST.ELA.K12.EE.3.1; CCSS.RI.4.1; TEKS.ELAR.4.7(C); CCSS.RI.5.1; TEKS.ELAR.5.7(C)] <br><br>3. Which event is described as a significant point of crisis during Kennedy's presidency? <br>A. The Bay of Pigs invasion <br>B. The Apollo 11 moon landing <br>C. The Cuban Missile Crisis <br>D. The Vietnam War escalation <br>[CCSS.RI.3.1; TEKS.ELAR.3.7(C); BEST.ELA.K12.EE.3.1; CCSS.RI.4.1; TEKS.ELAR.4.7(C); CCSS.RI.5.1; TEKS.ELAR.5.7(C)] <br><br>4. What action did Kennedy take in response to the Soviet missile bases found in Cuba? <br>A. He launched an immediate attack on Cuba <br>B. He blockaded Cuba&rsquo;s coast <br>C. He initiated peace talks with the Soviet Union <br>D. He ignored the situation entirely <br>[CCSS.RI.3.1; TEKS.ELAR.3.7(C); BEST.ELA.K12.EE.3.1; CCSS.RI.4.1; TEKS.ELAR.4.7(C); CCSS.RI.5.1; TEKS.ELAR.5.7(C)] <br><br>5. What was one consequence of America's actions during the Cold War, as indicated in the article? <br>A. Strengthened relations with Communist countries <br>B. Increased cooperation among NATO members <br>C. Escalation of military involvement in Vietnam <br>D. A decline in US technological advancements <br>[CCSS.RI.3.1; TEKS.ELAR.3.7(C); BEST.ELA.K12.EE.3.1; CCSS.RI.4.1; TEKS.ELAR.4.7(C); CCSS.RI.5.1; TEKS.ELAR.5.7(C)] <br><br></p>
'@

$s83 = @'
<p>In this informational e-book, Megan M. Gunderson describes the key plans and challenges faced by President John F. Kennedy during his administration through a detailed presentation of historical events and strategies. The text highlights Kennedy's initiatives regarding civil rights, the space race, and foreign relations amidst the Cold War context.</p>
'@

$s84 = @'
<ul>
<li>poverty (noun): the state of being extremely poor.</li>
<li>civil rights (noun): the rights of citizens to political and social freedom and equality.</li>
<li>cold war (proper noun): the state of political hostility between the Soviet bloc and the US-led Western powers.</li>
<li>soviet (proper adjective): relating to the former Soviet Union.</li>
<li>communist (proper adjective): relating to or advocating the political and economic system of communism.</li>
<li>space program (noun): a government-funded program aimed at the exploration of space.</li>
<li>relationships (noun): the way in which two or more people or groups regard and behave toward each other.</li>
<li>communism (proper noun): a political and economic ideology advocating for a classless system and collective ownership.</li>
<li>cubans (proper noun): citizens of Cuba.</li>
<li>bay of pigs (proper noun): a coastal area in Cuba where an unsuccessful invasion occurred in 1961.</li>
<li>failed (verb): did not succeed.</li>
<li>communist (proper adjective): relating to or advocating the political and economic system of communism.</li>
<li>india (proper noun): a country in South Asia.</li>
<li>vietnam (proper noun): a country in Southeast Asia.</li>
<li>advisers (noun): individuals who offer advice, especially in a professional capacity.</li>
<li>nuclear (adjective): relating to the nucleus of an atom or to the energy produced by nuclear reactions.</li>
<li>navy (proper noun): a branch of a nation's armed forces that conducts military operations at sea.</li>
<li>rights (noun): legal, social, or ethical principles of freedom or entitlement.</li>
<li>cuban missile crisis (proper noun): a 13-day confrontation between the United States and the Soviet Union in 1962 over the presence of missile sites in Cuba.</li>
</ul>
'@

$s85 = @'
<p>Making Plans (excerpt from John F. Kennedy) <br>By Megan M. Gunderson, Abdo Publishing <br><br>John F. Kennedy became President of the United States on January 20, 1961. He had many plans when he took office. These included programs for higher wages and aid for people experiencing *poverty*. He also hoped to improve *civil rights*. <br><br>During this time, America was part of the *Cold War* with the *Soviet* Union. Another plan was to put a person on the moon before the *Communist* Soviet Union did. As a result, Kennedy greatly improved the US *space program*. <br><br>*Relationships* with other countries were another concern for Kennedy. American officials wanted to stop the spread of *Communism*. So, they trained a group of *Cubans* to overthrow Cuba&rsquo;s *Communist* leader. In 1961, the group invaded Cuba at the *Bay of Pigs*. However, the effort *failed*. <br><br>Then, in 1962, *Communist* China invaded *India*. Kennedy sent weapons to India&rsquo;s army. Meanwhile, South *Vietnam* was fighting *Communist* North Vietnam. Kennedy sent thousands of US military *advisers* to help the South Vietnamese. <br><br>In October 1962, Kennedy learned *nuclear* missile bases were being built in Cuba. He wanted to stop *Soviet* ships from bringing supplies there. So, Kennedy stated that the US *Navy* would block Cuba&rsquo;s coast. <br><br>The Soviet leader said his country would guard its shipping *rights*. A *nuclear* war could begin. But Kennedy held his ground. The *Soviet* ships stayed away from the blocked area. Kennedy won the showdown! This event became known as the *Cuban Missile Crisis*.</p>
'@

$s86 = @'
<p>1. What was one of John F. Kennedy's major plans when he took office as President?<br>A. To dissolve the space program <br>B. To improve civil rights <br>C. To promote isolationism <br>D. To end the Cold War<br>[CCSS.RI.3.1; TEKS.ELAR.3.7(C)]<br><br>2. Which event depicted the efforts of the United States to intervene in Cuba during the Cold War?<br>A. The Cuban Missile Crisis <br>B. The Bay of Pigs invasion <br>C. The Space Race <br>D. The Vietnam War<br>[CCSS.RI.4.1; TEKS.ELAR.4.7(C)]<br><br>3. How did the Cold War affect U.S. military involvement in South Vietnam?<br>A. The U.S. decreased its military presence <br>B. The U.S. sent thousands of military advisers <br>C. The U.S. remained neutral <br>D. The U.S. declared war on China<br>[CCSS.RI.5.1; TEKS.ELAR.5.7(C)]<br><br>4. What was one result of Kennedy's response to the Soviet Union's actions in Cuba?<br>A. The Soviet Union retreated from Cuba <br>B. A nuclear war broke out <br>C. The Bay of Pigs invasion was successful <br>D. The U.S. signed a peace treaty with Cuba<br>[BEST.ELA.K12.EE.3.1; CCSS.RI.3.1]<br><br>5. What was the significance of the Cuban Missile Crisis during the Cold War?<br>A. It marked the end of the Cold War <br>B. It was a critical point in U.S.-Soviet relations <br>C. It led to an increase in space exploration funding <br>D. It resulted in the establishment of diplomatic relations with China<br>[CCSS.RI.4.1; TEKS.ELAR.4.7(C)]<br><br></p>
'@

$ws.Range("A32").Value2 = 'ARP_MakingPlans_G3-5'
$ws.Range("B32").Value2 = $s71
$ws.Range("C32").Value2 = 'John F. Kennedy; Cold War; Cuban Missile Crisis; US space program; civil rights'
$ws.Range("D32").Value2 = $s72
$ws.Range("E32").Value2 = $s73
$ws.Range("F32").Value2 = $s74

$ws.Range("A33").Value2 = 'ARP_MakingPlans_G3-5'
$ws.Range("B33").Value2 = $s71
$ws.Range("C33").Value2 = 'John F. Kennedy; Cold War; Cuban Missile Crisis; US space program; civil rights'
$ws.Range("D33").Value2 = $s72
$ws.Range("E33").Value2 = $s73
$ws.Range("F33").Value2 = $s74

$ws.Range("A34").Value2 = 'ARP_MakingPlans_G3-5'
$ws.Range("B34").Value2 = $s75
$ws.Range("C34").Value2 = 'John F. Kennedy; Cold War; Cuban Missile Crisis; space program; civil rights'
$ws.Range("D34").Value2 = $s76
$ws.Range("E34").Value2 = $s77
$ws.Range("F34").Value2 = $s78

$ws.Range("A35").Value2 = 'ARP_MakingPlans_G3-5'
$ws.Range("B35").Value2 = $s75
$ws.Range("C35").Value2 = 'John F. Kennedy; Cold War; Cuban Missile Crisis; space program; civil rights'
$ws.Range("D35").Value2 = $s76
$ws.Range("E35").Value2 = $s77
$ws.Range("F35").Value2 = $s78

$ws.Range("A36").Value2 = 'ARP_MakingPlans_G3-5'
$ws.Range("B36").Value2 = $s79
$ws.Range("C36").Value2 = 'John F. Kennedy; Cold War; Cuban Missile Crisis; civil rights; US space program'
$ws.Range("D36").Value2 = $s80
$ws.Range("E36").Value2 = $s81
$ws.Range("F36").Value2 = $s82

$ws.Range("A37").Value2 = 'ARP_MakingPlans_G3-5'
$ws.Range("B37").Value2 = $s79
$ws.Range("C37").Value2 = 'John F. Kennedy; Cold War; Cuban Missile Crisis; civil rights; US space program'
$ws.Range("D37").Value2 = $s80
$ws.Range("E37").Value2 = $s81
$ws.Range("F37").Value2 = $s82

$ws.Range("A38").Value2 = 'ARP_MakingPlans_G3-5'
$ws.Range("B38").Value2 = $s83
$ws.Range("C38").Value2 = 'John F. Kennedy; Cold War; Cuban Missile Crisis; space program; civil rights'
$ws.Range("D38").Value2 = $s84
$ws.Range("E38").Value2 = $s85
$ws.Range("F38").Value2 = $s86

$ws.Range("A39").Value2 = 'ARP_MakingPlans_G3-5'
$ws.Range("B39").Value2 = $s83
$ws.Range("C39").Value2 = 'John F. Kennedy; Cold War; Cuban Missile Crisis; space program; civil rights'
$ws.Range("D39").Value2 = $s84
$ws.Range("E39").Value2 = $s85
$ws.Range("F39").Value2 = $s86
